# Edit described by commit: "Modif cas d'utilisation + ajout table commande"
# - Adds a new "OrderDetail" table block (rows 62-69) to the data dictionary sheet
# - Re-merges rows 61 and 62 into full-width (B:Q) rows (blank separator + table title)
# - Updates the sheet view (scroll position / current selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-shape the merged cells for row 61 (blank separator) and row 62 (title)
#    Previously both rows were split into the usual 4 column-groups
#    (B:E / F:I / J:M / N:Q); now they each become one B:Q merge, matching
#    the pattern used by every other table separator/title in the sheet
#    (see rows 19/20, 28/29, 37/38/... -> B10:Q10, B20:Q20, etc.)
# ---------------------------------------------------------------------------
$ws.Range("B61:E61").UnMerge()
$ws.Range("F61:I61").UnMerge()
$ws.Range("J61:M61").UnMerge()
$ws.Range("N61:Q61").UnMerge()
$ws.Range("B61:Q61").Merge()

$ws.Range("B62:E62").UnMerge()
$ws.Range("F62:I62").UnMerge()
$ws.Range("J62:M62").UnMerge()
$ws.Range("N62:Q62").UnMerge()
$ws.Range("B62:Q62").Merge()

# ---------------------------------------------------------------------------
# 2. Fill in the new "OrderDetail" table definition (rows 62-69)
#    Row 62 is the table title (bold + centered, same style as the other
#    table headers such as B10 "Product", B20 "Supplier", B55 "Sailer" ...)
# ---------------------------------------------------------------------------
$ws.Range("B62").Value = "OrderDetail"
$ws.Range("B62").Font.Bold = $true
$ws.Range("B62").HorizontalAlignment = -4108

$ws.Range("B63").Value = "Numéro de commandes"
$ws.Range("F63").Value = "order_id"
$ws.Range("J63").Value = "Integer"
$ws.Range("N63").Value = ">0, Unique, AI"

$ws.Range("B64").Value = "Date de la commande"
$ws.Range("F64").Value = "order_date"
$ws.Range("J64").Value = "Date"

$ws.Range("B65").Value = "Nombre de produit commandé "
$ws.Range("F65").Value = "product_amount"
$ws.Range("J65").Value = "Integer"

$ws.Range("B66").Value = "Sous total HT par produit"
$ws.Range("F66").Value = "product_sub_total"
$ws.Range("J66").Value = "Numeric"

$ws.Range("B67").Value = "Total HT"
$ws.Range("F67").Value = "total_ht"
$ws.Range("J67").Value = "Numeric"

$ws.Range("B68").Value = "Total TTC"
$ws.Range("F68").Value = "total_ttc"
$ws.Range("J68").Value = "Numeric"

$ws.Range("B69").Value = "Total des taxes"
$ws.Range("F69").Value = "total_tva"
$ws.Range("J69").Value = "Numeric"

# The "N" column formulas/notes are filled in last, in this exact order,
# to reproduce the shared-string creation order observed in the target
# workbook (product_amount * pro_price = 106, Sum(product_sub_total) = 107,
# total_ht * tva = 108, total_ht * (1+tva) = 109)
$ws.Range("N66").Value = "product_amount * pro_price"
$ws.Range("N67").Value = "Sum(product_sub_total)"
$ws.Range("N69").Value = "total_ht * tva"
$ws.Range("N68").Value = "total_ht * (1+tva)"

# ---------------------------------------------------------------------------
# 3. Update the sheet view: scroll to row 40 and select N70:Q70
# ---------------------------------------------------------------------------
$ws.Range("N70:Q70").Select()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
